$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest decade of data (2000年-2009年, rows 2-11). Deleting the
# entire rows shifts the following rows (2010年-2018年, old rows 12-20) up
# so they become the new rows 2-10.
$ws.Range("A2:E11").EntireRow.Delete()
